# The author switched from the "Lecturer Free" tab back to the
# "ClassroomA" tab and selected the classroom's available-times column
# (A2:B67, i.e. everything below the header row) while continuing work
# on limiting classroom/lecturer available times.

$wb = $excel.ActiveWorkbook

$wsClassroomA = $wb.Worksheets.Item("ClassroomA")
$wsClassroomA.Activate() | Out-Null
$wsClassroomA.Range("A2:B67").Select() | Out-Null
